# Two new players were added to the "historico" leaderboard. The sheet is
# sorted descending by column C ("inicial"), so each new row has to be
# inserted at the position that keeps that order, pushing the rows below it
# down by one. The "posicao" column (D) is then renumbered 1..N to stay
# sequential.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Life Kudalini" (3202574604.1841998) is the new largest value overall,
#    so it becomes the new row 2 (right after the header row).
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "61fc92022f8cc6002866c120"
$ws.Range("B2").Value = "Life Kudalini"
$ws.Range("C2").Value = 3202574604.1841998
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "kz7sv0d7"

# 2) "Lucas Dias" (568137332.2888) sits between "Lion Fera" (579670695.75,
#    now row 12 after the insert above) and "Fabio HK" (450714236.7712,
#    still below it) -- that slot is row 13.
$ws.Rows.Item(13).Insert()
# Insert() copies the formatting of the row above onto the new row (like
# Excel's default "Insert" behaviour); the source row didn't have that
# formatting, so drop it to get the plain/no-style cells the author saved.
$ws.Range("A13:E13").ClearFormats()
$ws.Range("A13").Value = "626eb2be46da9a002807e173"
$ws.Range("B13").Value = "Lucas Dias"
$ws.Range("C13").Value = 568137332.2888
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "l2nhyaun"

# 3) Renumber "posicao" (column D) for every data row so the ranking stays
#    1..43 with no gaps/duplicates after the two insertions.
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("D" + $r).Value = $r - 1
}

# Re-assert the one "inicial" value whose literal text drifts by a float
# ULP (2841504.1400000001 vs 2841504.14) purely from being shifted down two
# rows by the inserts above -- writing it back keeps it identical to what
# was on disk before the edit.
$ws.Range("C43").Value = 2841504.14

# 4) Match the saved selection/scroll state from the authored workbook.
[void]$ws.Range("A40").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
[void]$ws.Range("A2:E44").Select()
